$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.669.97'
$ws.Range("E2").Value = '  +1.69%  '

# Row 3
$ws.Range("D3").Value = '2.302.73'
$ws.Range("E3").Value = '  +0.46%  '

# Row 4
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.67'
$ws.Range("E5").Value = '  -0.08%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.79'
$ws.Range("E6").Value = '  +0.50%  '

# Row 7
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("E8").Value = '  +0.15%  '

# Row 9
$ws.Range("E9").Value = '  +0.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.98'
$ws.Range("E10").Value = '  +1.78%  '

# Row 11
$ws.Range("E11").Value = '  +0.01%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.55'
$ws.Range("E12").Value = '  +3.52%  '

# Row 13
$ws.Range("E13").Value = '  +0.79%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.998'
$ws.Range("E14").Value = '  +4.06%  '

# Row 15
$ws.Range("E15").Value = '  +0.84%  '

# Row 16
$ws.Range("D16").Value = '2.652.75'
$ws.Range("E16").Value = '  +0.51%  '

# Row 17
$ws.Range("D17").Value = '2.305.67'
$ws.Range("E17").Value = '  +0.47%  '

# Row 18
$ws.Range("D18").Value = '42.491.40'
$ws.Range("E18").Value = '  +1.26%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.59'
$ws.Range("E19").Value = '  +3.00%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000105'
$ws.Range("E20").Value = '  -0.17%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.66'
$ws.Range("E21").Value = '  +33.96%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.98'
$ws.Range("E22").Value = '  +1.00%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.55'
$ws.Range("E23").Value = '  -2.11%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '267.40'
$ws.Range("E24").Value = '  -3.98%  '

# Row 25
$ws.Range("E25").Value = '  -0.77%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.33%  '

# Row 27
$ws.Range("E27").Value = '  +1.33%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.26'
$ws.Range("E28").Value = '  -6.60%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.58'
$ws.Range("E29").Value = '  -0.90%  '

# Row 30
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.64'
$ws.Range("E30").Value = '  +14.56%  '

# Row 31
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.03'
$ws.Range("E31").Value = '  +5.84%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.35'
$ws.Range("E32").Value = '  +1.35%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0882'
$ws.Range("E33").Value = '  +1.34%  '

# Row 34
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.67'
$ws.Range("E34").Value = '  -6.14%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.132'
$ws.Range("E35").Value = '  -3.33%  '

# Row 36
$ws.Range("E36").Value = '  +0.28%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.59'
$ws.Range("E37").Value = '  +1.97%  '

# Row 39
$ws.Range("E39").Value = '  -0.43%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.76'
$ws.Range("E40").Value = '  -2.21%  '

# Row 41
$ws.Range("E41").Value = '  +13.73%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.07'
$ws.Range("E42").Value = '  -0.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.27'
$ws.Range("E43").Value = '  +1.28%  '

# Row 44
$ws.Range("E44").Value = '  +0.76%  '

# Row 45
$ws.Range("E45").Value = '  -0.06%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '117.17'
$ws.Range("E46").Value = '  +4.23%  '

# Row 47
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.35'
$ws.Range("E47").Value = '  +3.80%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '80.30'
$ws.Range("E48").Value = '  +4.31%  '

# Row 49
$ws.Range("D49").Value = '1.639.08'
$ws.Range("E49").Value = '  +3.53%  '

# Row 50
$ws.Range("E50").Value = '  +0.48%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.87'
$ws.Range("E51").Value = '  -0.24%  '
